# Updates the "Översikt VANSBRO" log sheet:
#  1. Column C ("Förändrad") timestamp bumped from 45184 to 45186 for every
#     data row (rows 2-297).
#  2. Each HYPERLINK() formula in columns S:Y that only has the URL
#     argument gets a second "friendly name" argument equal to the row's
#     "Beteckning" (column A) value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count - 1   # sheet's row 1 (XML row index 0) is a blank spacer row above the header

for ($row = 2; $row -le $lastRow; $row++) {

    # 1) Bump the "Förändrad" date in column C (45184 -> 45186) when present.
    $cCell = $ws.Cells.Item($row, 3)
    $cVal = $cCell.Value2
    if ($cVal -eq 45184) {
        $cCell.Value = 45186
    }

    # 2) Add the friendly-name second argument to HYPERLINK formulas (cols S:Y = 19:25).
    $name = $ws.Cells.Item($row, 1).Value2
    if ($name) {
        for ($col = 19; $col -le 25; $col++) {
            $cell = $ws.Cells.Item($row, $col)
            $formula = $cell.Formula
            if ($formula -and $formula -match '^=HYPERLINK\("([^"]*)"\)$') {
                $url = $matches[1]
                $cell.Formula = '=HYPERLINK("' + $url + '", "' + $name + '")'
            }
        }
    }
}
